# Weekly update: a new price-observation date (2023-01-08, serial 44931) was
# added for the "Terminal La Palmera de La Serena" / Mango series. It is
# inserted as a new group of 3 rows (Especial/Primera/Segunda) right above
# the existing 2022-03-15 (44545) group, i.e. at row 931, pushing every row
# from 931 downward by 3. The new group repeats the most recent prior
# observation's prices/volume/origin (the 44613 group, Perú) under the new
# date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 931 (shifts rows 931:1011 down to 934:1014 and
# the sheet dimension grows from A1:T1011 to A1:T1014 automatically).
$ws.Range("A931:A933").EntireRow.Insert()

# After the insert, the data that used to sit at rows 943:945 (date 44613,
# Perú, Especial/Primera/Segunda) now lives at rows 946:948. Clone those
# three rows into the freshly inserted 931:933, then overwrite the date
# column with the new observation date (44931).
for ($i = 0; $i -lt 3; $i++) {
    $srcRow = 946 + $i
    $dstRow = 931 + $i
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($dstRow, $col).Value2 = $ws.Cells.Item($srcRow, $col).Value2
    }
    # Column D = Fecha (date serial)
    $ws.Cells.Item($dstRow, 4).Value2 = 44931
}
